# A2osX Issue List - update script
# - Rename "Sheet1" -> "Memory"
# - Switch active tab from Issues to Suggestions
# - Add new Issue row (Issue #176) to the Issues sheet
# - Update the "Rename A2OSX.BUILD..." suggestion text on the Suggestions sheet
# - Update selections to match the saved workbook state

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Memory ---
$wsMemory = $wb.Worksheets.Item("Sheet1")
$wsMemory.Name = "Memory"

# --- Issues sheet: add new row 85 ---
$wsIssues = $wb.Worksheets.Item("Issues")
$wsIssues.Range("B85").Value = 176
$wsIssues.Range("C85").Value = "EDIT"
$wsIssues.Range("D85").Value = 1129
$wsIssues.Range("E85").Value = "If are editing a new file ""EDIT NEWFILE"" when you go to save the first time the name field is blank, you have to enter the name again.  Update: this happens multiple times while editing the file."
$wsIssues.Rows.Item(85).RowHeight = 75

# --- Suggestions sheet: update text of row 65 (A2OSX.BUILD rename suggestion) ---
$wsSuggestions = $wb.Worksheets.Item("Suggestions")
$wsSuggestions.Range("D65").Value = "Rename A2OSX.BUILD to A2DEV, it saves memory and shortens PATH.  Also A2DEV reflects that it is a development image, not a production/test/release image.  Move files in BIN/DEV to BIN and remove BIN/DEV from PATH (done in INIT).  Add USR/BIN to PATH (for User Scripts and programs)."
$wsSuggestions.Rows.Item(65).RowHeight = 90

# --- Selection / active sheet state ---
$wsIssues.Activate()
$wsIssues.Range("E86").Select()

$wsSuggestions.Activate()
$wsSuggestions.Range("D66").Select()
